$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'326.36"
$ws.Range("E2").Value = "'-0.54%"

$ws.Range("D3").Value = "'43.72"
$ws.Range("E3").Value = "'-0.38%"

$ws.Range("D4").Value = "'5.553"
$ws.Range("E4").Value = "'-0.17%"

$ws.Range("D5").Value = "'0.07999"
$ws.Range("E5").Value = "'-1.19%"

$ws.Range("D6").Value = "'1.940"
$ws.Range("E6").Value = "'1.19%"

$ws.Range("B7").Value = "BTSEToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D7").Value = "'2.543"
$ws.Range("E7").Value = "'-7.80%"

$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9420"
$ws.Range("E8").Value = "'-0.69%"

$ws.Range("D9").Value = "'0.1151"
$ws.Range("E9").Value = "'-3.01%"

$ws.Range("D10").Value = "'0.1837"
$ws.Range("E10").Value = "'-3.35%"

$ws.Range("D11").Value = "'11.99"
$ws.Range("E11").Value = "'37.78%"

$ws.Range("D12").Value = "'0.09654"
$ws.Range("E12").Value = "'-0.34%"

$ws.Range("D13").Value = "'0.04755"
$ws.Range("E13").Value = "'16.03%"

$ws.Range("E14").Value = "'-0.06%"

$ws.Range("E15").Value = "'-1.03%"

$ws.Range("D16").Value = "'0.04071"
$ws.Range("E16").Value = "'-5.79%"

$ws.Range("D17").Value = "'0.005986"
$ws.Range("E17").Value = "'1.02%"

$ws.Range("D18").Value = "'3.376"
$ws.Range("E18").Value = "'-5.68%"

$ws.Range("D19").Value = "'4.318"
$ws.Range("E19").Value = "'0.57%"

$ws.Range("D20").Value = "'0.3474"

$ws.Range("D21").Value = "'0.1404"
$ws.Range("E21").Value = "'2.88%"

$ws.Range("E22").Value = "'-2.82%"

$ws.Range("D23").Value = "'0.001250"
$ws.Range("E23").Value = "'0.72%"

$ws.Range("D24").Value = "'0.004289"
$ws.Range("E24").Value = "'-7.83%"

$ws.Range("D25").Value = "'0.0001195"
$ws.Range("E25").Value = "'-3.07%"

$ws.Range("D26").Value = "'0.0003764"

$ws.Range("D38").Value = "'0.02543"
$ws.Range("E38").Value = "'-4.59%"

$ws.Range("D39").Value = "'0.05426"
$ws.Range("E39").Value = "'-1.11%"

$ws.Range("D40").Value = "'0.007508"
$ws.Range("E40").Value = "'-2.74%"

$ws.Range("D41").Value = "'0.1384"
$ws.Range("E41").Value = "'-0.87%"

$ws.Range("D42").Value = "'0.007480"
$ws.Range("E42").Value = "'-34.01%"

$ws.Range("D43").Value = "'0.002025"
$ws.Range("E43").Value = "'-4.18%"

$ws.Range("D44").Value = "'0.008319"
$ws.Range("E44").Value = "'-11.93%"

$ws.Range("D45").Value = "'0.00007127"
$ws.Range("E45").Value = "'1.46%"

$ws.Range("E46").Value = "'0.28%"

$ws.Range("E47").Value = "'1.22%"

$ws.Range("D48").Value = "'0.003474"
$ws.Range("E48").Value = "'0.73%"

$ws.Range("E49").Value = "'0.28%"

$ws.Range("D50").Value = "'0.0002011"
$ws.Range("E50").Value = "'0.28%"
